# Fruta / hortaliza, semanal
# Update the weekly price records: the D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado) and S (Precio $/Kg) values
# for rows 2-8 are reshuffled to reflect the new weekly data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row -> values (D, M, N, O, P, S)
$rows = @{
    2 = @(44253, 90, 12000, 13000, 12667, 905)
    3 = @(44172, 90, 8500,  9000,  8806,  629)
    4 = @(44232, 60, 11000, 12000, 11583, 827)
    5 = @(44216, 55, 11000, 12000, 11545, 825)
    6 = @(44181, 65, 9000,  10000, 9462,  676)
    7 = @(44210, 70, 10000, 11000, 10357, 740)
    8 = @(44229, 55, 11000, 12000, 11364, 812)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 13).Value = $vals[1]  # M - Volumen
    $ws.Cells.Item($r, 14).Value = $vals[2]  # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals[3]  # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals[4]  # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $vals[5]  # S - Precio $/Kg
}
